$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# The "water_level_ini_type" column (L) is being removed — the default
# initial water level agg type will now be read from sqlite ("min" instead
# of "max"), so the column is no longer needed here. Delete the entire
# column; everything to the right shifts left by one.
$ws.Columns("L").Delete()

# Update the active selection / scroll position left behind by the edit
# (was topLeftCell=Z1 / AJ12 selected; now back near the top-left, K5).
$ws.Range("K5").Select()
